# Applies the log-content edits described by the commit:
#   "add animation for chatbot, make new error log for windows error"
#
# The sheet is a flat operation/error log. The author:
#   - renamed the user from "Haruka Murakami" to "Momoko Shimizu" on every row
#   - renumbered/renamed the screenshot file names in column J
#     (bdotYYYYMMDD_HHMMSS/...png) to a new capture session with simple
#     sequential file names
#   - rewrote several of the Japanese explanation texts in column K
#   - moved the Windows-Update-related "error" row earlier in the sequence
#     (row 5 is now the error row, row 7 is now back to a normal operation
#     row) and gave the new error row (row 5) fresh error_type/error_content
#     values in columns L/M, clearing out the old error row's (row 7) L/M
#     values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: user_name, every data row (2-16) becomes "Momoko Shimizu" ---
foreach ($r in 2..16) {
    $ws.Cells.Item($r, 3).Value = "Momoko Shimizu"
}

# --- Row 2 ---
$ws.Range("J2").Value = "bdot20240415_141954/1.png"
$ws.Range("K2").Value = "「スタート」ボタンをクリックする"

# --- Row 3 ---
$ws.Range("J3").Value = "bdot20240415_141954/2.png"
$ws.Range("K3").Value = "メニューから「設定」アイコンをクリックする"

# --- Row 4 ---
$ws.Range("J4").Value = "bdot20240415_141954/3.png"
$ws.Range("K4").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"

# --- Row 5: becomes the new error row ---
$ws.Range("B5").Value = "error"
$ws.Range("J5").Value = "bdot20240415_141954/4.png"
$ws.Range("K5").Value = "0x80240fff エラー"
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"

# --- Row 6 ---
$ws.Range("J6").Value = "bdot20240415_141954/5.png"
$ws.Range("K6").Value = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"

# --- Row 7: was the error row, now a normal operation row again ---
$ws.Range("B7").Value = "operation"
$ws.Range("J7").Value = "bdot20240415_141954/5.png"
$ws.Range("K7").Value = "メニューからターミナル(管理者)をクリックする"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""

# --- Row 8 ---
$ws.Range("J8").Value = "bdot20240415_141954/6.png"
$ws.Range("K8").Value = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"

# --- Row 9 ---
$ws.Range("J9").Value = "bdot20240415_141954/7.png"
$ws.Range("K9").Value = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"

# --- Row 10 ---
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"

# --- Row 11 ---
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"

# --- Row 12 ---
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"

# --- Row 13 ---
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"

# --- Row 14 ---
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"

# --- Row 15 ---
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"

# --- Row 16 ---
$ws.Range("J16").Value = "bdot20240415_141954/11.png"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"
